$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control_panel")

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 6
$ws.Range("H10").Value = 7
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 9
